# Refresh the clan/members roster (Name, Platform) pulled from the
# Selenium-based scraper. The sheet already has the header row
# (Name/Platform) in row 1, so we just rewrite the member rows below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$members = @(
    @("Alexvomits", "epic"),
    @("Anubisgoat1", "epic"),
    @("Army gunners", "epic"),
    @("AVN7", "epic"),
    @("BeerMasjien", "epic"),
    @("Bobman2947", "xbl"),
    @("BRO_Dragaawn", "epic"),
    @("Chutney_Tjops", "epic"),
    @("Creativ1s", "epic"),
    @("daCostaRaps", "psn"),
    @("Eskapa1d", "epic"),
    @("FiNniCKin_2Jz", "psn"),
    @("GreenDragon0603", "epic"),
    @("Hannib8l OG", "epic"),
    @("HerezAJoHnnY", "epic"),
    @("Jamalvandux", "epic"),
    @("Japes360", "psn"),
    @("K1LL1NGF13LD5", "epic"),
    @("KillerPranesh", "epic"),
    @("MightyMidget8811", "psn"),
    @("MonsterGames200", "epic"),
    @("Mr_Ballistic_", "epic"),
    @("Nanashi_ZA", "epic"),
    @("NaVeDs", "epic"),
    @("NoNoob1415", "epic"),
    @("Oops123rage", "xbl"),
    @("Reckless_ness11", "epic"),
    @("Rusting", "epic"),
    @("SONIC-_1397", "epic"),
    @("sweetdreamcringe", "epic"),
    @("TakenPersonally", "epic"),
    @("THAPZIN_6", "epic"),
    @("TheBromeo", "xbl"),
    @("Turtlejuice 2.0", "epic"),
    @("Voorie", "epic"),
    @("Wasabi_ZA", "epic"),
    @("ZeroTwo0_ ", "epic")
)

$row = 2
foreach ($member in $members) {
    $ws.Cells.Item($row, 1).Value = $member[0]
    $ws.Cells.Item($row, 2).Value = $member[1]
    $row = $row + 1
}
